$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 30   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  5/1/2023  Through  5/7/2023"

# --- Row 30 (Hate Crimes): D/E become text placeholders ("0" / "***.*") ---
$ws.Range("D30").Value = "'0"
$ws.Range("E30").Value = "***.*"
$ws.Range("M26").Copy()
$ws.Range("D30:E30").PasteSpecial(-4122)

# --- Numeric data updates, rows 14-30 ---
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = -50
$ws.Range("G14").Value = 6
$ws.Range("H14").Value = -16.666666666666
$ws.Range("I14").Value = 22
$ws.Range("J14").Value = 20
$ws.Range("K14").Value = 10
$ws.Range("L14").Value = -18.518518518518
$ws.Range("M14").Value = -45
$ws.Range("N14").Value = -86.075949367088

$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = -80
$ws.Range("F15").Value = 15
$ws.Range("G15").Value = 13
$ws.Range("H15").Value = 15.384615384615
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = -2.5
$ws.Range("L15").Value = 6.849315068493
$ws.Range("M15").Value = -3.703703703703
$ws.Range("N15").Value = -63.551401869158

$ws.Range("C16").Value = 52
$ws.Range("D16").Value = 65
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 193
$ws.Range("G16").Value = 193
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 803
$ws.Range("J16").Value = 835
$ws.Range("K16").Value = -3.832335329341
$ws.Range("L16").Value = 30.357142857142
$ws.Range("M16").Value = -26.465201465201
$ws.Range("N16").Value = -85.65815324165

$ws.Range("C17").Value = 85
$ws.Range("D17").Value = 78
$ws.Range("E17").Value = 8.974358974358
$ws.Range("F17").Value = 326
$ws.Range("G17").Value = 299
$ws.Range("H17").Value = 9.030100334448
$ws.Range("I17").Value = 1349
$ws.Range("J17").Value = 1304
$ws.Range("K17").Value = 3.450920245398
$ws.Range("L17").Value = 30.464216634429
$ws.Range("M17").Value = 27.50472589792
$ws.Range("N17").Value = -50.440852314474

$ws.Range("C18").Value = 41
$ws.Range("E18").Value = -26.785714285714
$ws.Range("F18").Value = 146
$ws.Range("G18").Value = 182
$ws.Range("H18").Value = -19.780219780219
$ws.Range("I18").Value = 720
$ws.Range("J18").Value = 847
$ws.Range("K18").Value = -14.994096812278
$ws.Range("L18").Value = 14.649681528662
$ws.Range("M18").Value = -22.162162162162
$ws.Range("N18").Value = -81.547924141465

$ws.Range("C19").Value = 104
$ws.Range("D19").Value = 119
$ws.Range("E19").Value = -12.605042016806
$ws.Range("F19").Value = 441
$ws.Range("G19").Value = 447
$ws.Range("H19").Value = -1.342281879194
$ws.Range("I19").Value = 1938
$ws.Range("J19").Value = 1877
$ws.Range("K19").Value = 3.249866808737
$ws.Range("L19").Value = 39.424460431654
$ws.Range("M19").Value = 49.537037037037
$ws.Range("N19").Value = -9.14205344585

$ws.Range("C20").Value = 29
$ws.Range("D20").Value = 27
$ws.Range("E20").Value = 7.407407407407
$ws.Range("G20").Value = 114
$ws.Range("H20").Value = 17.543859649122
$ws.Range("I20").Value = 565
$ws.Range("J20").Value = 593
$ws.Range("K20").Value = -4.721753794266
$ws.Range("L20").Value = 34.204275534441
$ws.Range("M20").Value = 24.449339207048
$ws.Range("N20").Value = -83.387239047339

$ws.Range("C21").Value = 313
$ws.Range("D21").Value = 352
$ws.Range("E21").Value = -11.079545454545
$ws.Range("F21").Value = 1260
$ws.Range("G21").Value = 1254
$ws.Range("H21").Value = 0.478468899521
$ws.Range("I21").Value = 5475
$ws.Range("J21").Value = 5556
$ws.Range("K21").Value = -1.45788336933
$ws.Range("L21").Value = 30.699450942945
$ws.Range("M21").Value = 10.695511524464
$ws.Range("N21").Value = -69.799768326989

$ws.Range("C22").Value = 6
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = 50
$ws.Range("F22").Value = 25
$ws.Range("G22").Value = 25
$ws.Range("I22").Value = 109
$ws.Range("J22").Value = 131
$ws.Range("K22").Value = -16.793893129771
$ws.Range("L22").Value = 31.325301204819
$ws.Range("M22").Value = -21.582733812949

$ws.Range("C23").Value = 23
$ws.Range("D23").Value = 35
$ws.Range("E23").Value = -34.285714285714
$ws.Range("F23").Value = 115
$ws.Range("G23").Value = 112
$ws.Range("H23").Value = 2.678571428571
$ws.Range("I23").Value = 546
$ws.Range("J23").Value = 493
$ws.Range("K23").Value = 10.750507099391
$ws.Range("L23").Value = 19.474835886214
$ws.Range("M23").Value = 60.117302052785

$ws.Range("C24").Value = 253
$ws.Range("D24").Value = 244
$ws.Range("E24").Value = 3.688524590163
$ws.Range("F24").Value = 948
$ws.Range("G24").Value = 1008
$ws.Range("H24").Value = -5.95238095238
$ws.Range("I24").Value = 4186
$ws.Range("J24").Value = 4269
$ws.Range("K24").Value = -1.944249238697
$ws.Range("L24").Value = 27.583053946967
$ws.Range("M24").Value = 29.838709677419

$ws.Range("C25").Value = 133
$ws.Range("D25").Value = 102
$ws.Range("E25").Value = 30.392156862745
$ws.Range("F25").Value = 461
$ws.Range("G25").Value = 469
$ws.Range("H25").Value = -1.705756929637
$ws.Range("I25").Value = 2005
$ws.Range("J25").Value = 2017
$ws.Range("K25").Value = -0.59494298463
$ws.Range("L25").Value = 44.140905823148
$ws.Range("M25").Value = -22.347017815646

$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -62.5
$ws.Range("F26").Value = 20
$ws.Range("G26").Value = 20
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 111
$ws.Range("J26").Value = 123
$ws.Range("K26").Value = -9.756097560975
$ws.Range("L26").Value = -11.2

$ws.Range("C27").Value = 13
$ws.Range("D27").Value = 6
$ws.Range("E27").Value = 116.666666666667
$ws.Range("F27").Value = 57
$ws.Range("G27").Value = 47
$ws.Range("H27").Value = 21.27659574468
$ws.Range("I27").Value = 208
$ws.Range("J27").Value = 197
$ws.Range("K27").Value = 5.583756345177
$ws.Range("L27").Value = -6.306306306306

$ws.Range("C28").Value = 8
$ws.Range("D28").Value = 8
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 18
$ws.Range("H28").Value = -35.714285714285
$ws.Range("I28").Value = 78
$ws.Range("J28").Value = 87
$ws.Range("K28").Value = -10.344827586206
$ws.Range("L28").Value = -32.173913043478
$ws.Range("M28").Value = -47.297297297297
$ws.Range("N28").Value = -88.055130168453

$ws.Range("C29").Value = 6
$ws.Range("D29").Value = 6
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 15
$ws.Range("H29").Value = -34.782608695652
$ws.Range("I29").Value = 66
$ws.Range("J29").Value = 77
$ws.Range("K29").Value = -14.285714285714
$ws.Range("L29").Value = -34
$ws.Range("M29").Value = -41.071428571428
$ws.Range("N29").Value = -88.888888888888

$ws.Range("F30").Value = 5
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 150
$ws.Range("I30").Value = 29
$ws.Range("K30").Value = 26.086956521739
$ws.Range("L30").Value = 61.111111111111
